$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 850
$ws.Range("I18").Value = 850
$ws.Range("K18").Value = 850
$ws.Range("M18").Value = -566

$ws.Range("H64").Value = 45457536
$ws.Range("I64").Value = 111113144
$ws.Range("J64").Value = 3650.7693
$ws.Range("K64").Value = 111113144
$ws.Range("L64").Value = 3650.7693
$ws.Range("M64").Value = -111112896
$ws.Range("N64").Value = -4146.7693

$ws.Range("H67").Value = 45457536
$ws.Range("I67").Value = 111113144
$ws.Range("J67").Value = 3650.7693
$ws.Range("K67").Value = 111113144
$ws.Range("L67").Value = 3650.7693
$ws.Range("M67").Value = -111112286
$ws.Range("N67").Value = -5366.7693

$ws.Range("H69").Value = 3890.1
$ws.Range("I69").Value = 3828.25
$ws.Range("J69").Value = 3896.9722
$ws.Range("K69").Value = 11484.75
$ws.Range("L69").Value = 11690.9166
$ws.Range("M69").Value = -10610.75
$ws.Range("N69").Value = -13438.9166

$ws.Range("H72").Value = 3890.1
$ws.Range("I72").Value = 3828.25
$ws.Range("J72").Value = 3896.9722
$ws.Range("K72").Value = 34454.25
$ws.Range("L72").Value = 35072.74980000001
$ws.Range("M72").Value = -30086.25
$ws.Range("N72").Value = -43808.74980000001

$ws.Range("H74").Value = 3190.1365
$ws.Range("I74").Value = 2535.375
$ws.Range("J74").Value = 3564.2856
$ws.Range("K74").Value = 2535.375
$ws.Range("L74").Value = 3564.2856
$ws.Range("M74").Value = -1599.375
$ws.Range("N74").Value = -5436.2856

$ws.Range("H77").Value = 3190.1365
$ws.Range("I77").Value = 2535.375
$ws.Range("J77").Value = 3564.2856
$ws.Range("K77").Value = 12676.875
$ws.Range("L77").Value = 17821.428
$ws.Range("M77").Value = -7996.875
$ws.Range("N77").Value = -27181.428

$ws.Range("H96").Value = 2050
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 1100
$ws.Range("K96").Value = 9000
$ws.Range("L96").Value = 3300
$ws.Range("M96").Value = -7627
$ws.Range("N96").Value = -6046

$ws.Range("H100").Value = 3061
$ws.Range("I100").Value = 1846.125
$ws.Range("K100").Value = 1846.125
$ws.Range("M100").Value = -1305.125

$ws.Range("H111").Value = 917.6
$ws.Range("I111").Value = 897
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 2691
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 376
$ws.Range("N111").Value = -9134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 765.1539
$ws.Range("I97").Value = 662.9286
$ws.Range("J97").Value = 884.4167
$ws.Range("K97").Value = 662.9286
$ws.Range("L97").Value = 884.4167
$ws.Range("M97").Value = -166.9286
$ws.Range("N97").Value = -1876.4167

$ws.Range("H102").Value = 1107.5
$ws.Range("I102").Value = 1107.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1107.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 514.5
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1074.25
$ws.Range("I122").Value = 1062.8182
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3188.4546
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -738.4546
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 971.5925999999999
$ws.Range("I94").Value = 775.2222
$ws.Range("J94").Value = 1364.3334
$ws.Range("K94").Value = 775.2222
$ws.Range("L94").Value = 1364.3334
$ws.Range("M94").Value = -324.2222
$ws.Range("N94").Value = -2266.3334

$ws.Range("H105").Value = 2027.2142
$ws.Range("I105").Value = 1241.6666
$ws.Range("J105").Value = 2616.375
$ws.Range("K105").Value = 1241.6666
$ws.Range("L105").Value = 2616.375
$ws.Range("M105").Value = 505.3334
$ws.Range("N105").Value = -6110.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3594.3333
$ws.Range("I134").Value = 3633.7778
$ws.Range("J134").Value = 3476
$ws.Range("K134").Value = 10901.3334
$ws.Range("L134").Value = 10428
$ws.Range("M134").Value = -8366.3334
$ws.Range("N134").Value = -15498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 684.5454999999999
$ws.Range("I117").Value = 182.5
$ws.Range("J117").Value = 971.4286
$ws.Range("K117").Value = 547.5
$ws.Range("L117").Value = 2914.2858
$ws.Range("M117").Value = 2894.5
$ws.Range("N117").Value = -9798.2858

$ws.Range("H122").Value = 1193.7273
$ws.Range("J122").Value = 2449.25
$ws.Range("L122").Value = 22043.25
$ws.Range("N122").Value = -26943.25

$ws.Range("H140").Value = 16608.658
$ws.Range("I140").Value = 31464.324
$ws.Range("J140").Value = 2514.8206
$ws.Range("K140").Value = 94392.97200000001
$ws.Range("L140").Value = 7544.4618
$ws.Range("M140").Value = -89212.97200000001
$ws.Range("N140").Value = -17904.4618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 720.8461
$ws.Range("I97").Value = 720.8461
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 720.8461
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -224.8461
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 6395
$ws.Range("I102").Value = 10898.4
$ws.Range("J102").Value = 3178.2856
$ws.Range("K102").Value = 10898.4
$ws.Range("L102").Value = 3178.2856
$ws.Range("M102").Value = -9276.4
$ws.Range("N102").Value = -6422.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2682.8572
$ws.Range("I7").Value = 2601.4285
$ws.Range("J7").Value = 2764.2856
$ws.Range("K7").Value = 2601.4285
$ws.Range("L7").Value = 2764.2856
$ws.Range("M7").Value = -2489.4285
$ws.Range("N7").Value = -2988.2856

$ws.Range("H40").Value = 1725.1613
$ws.Range("I40").Value = 1388.7391
$ws.Range("J40").Value = 2692.375
$ws.Range("K40").Value = 1388.7391
$ws.Range("L40").Value = 2692.375
$ws.Range("M40").Value = -1252.7391
$ws.Range("N40").Value = -2964.375

$ws.Range("H61").Value = 2173.9375
$ws.Range("I61").Value = 1444.7778
$ws.Range("J61").Value = 3111.4285
$ws.Range("K61").Value = 1444.7778
$ws.Range("L61").Value = 3111.4285
$ws.Range("M61").Value = -1242.7778
$ws.Range("N61").Value = -3515.4285

$ws.Range("H93").Value = 1612.5
$ws.Range("I93").Value = 1150
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1150
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = 98
$ws.Range("N93").Value = -5496

$ws.Range("H100").Value = 3125.7144
$ws.Range("I100").Value = 3040
$ws.Range("J100").Value = 3340
$ws.Range("K100").Value = 3040
$ws.Range("L100").Value = 3340
$ws.Range("M100").Value = -2499
$ws.Range("N100").Value = -4422

$ws.Range("H113").Value = 2173.9375
$ws.Range("I113").Value = 1444.7778
$ws.Range("J113").Value = 3111.4285
$ws.Range("K113").Value = 1444.7778
$ws.Range("L113").Value = 3111.4285
$ws.Range("M113").Value = 725.2221999999999
$ws.Range("N113").Value = -7451.4285

$ws.Range("H121").Value = 35140
$ws.Range("J121").Value = 35140
$ws.Range("L121").Value = 35140
$ws.Range("N121").Value = -38634

$ws.Range("H122").Value = 53085.05
$ws.Range("I122").Value = 74085.78999999999
$ws.Range("J122").Value = 4083.3333
$ws.Range("K122").Value = 222257.37
$ws.Range("L122").Value = 12249.9999
$ws.Range("M122").Value = -219807.37
$ws.Range("N122").Value = -17149.9999

$ws.Range("H126").Value = 2682.8572
$ws.Range("I126").Value = 2601.4285
$ws.Range("J126").Value = 2764.2856
$ws.Range("K126").Value = 7804.2855
$ws.Range("L126").Value = 8292.856800000001
$ws.Range("M126").Value = -5334.2855
$ws.Range("N126").Value = -13232.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2339.5
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 2619.3333
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2619.3333
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -5365.3333

$ws.Range("H107").Value = 423.125
$ws.Range("I107").Value = 426.42856
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1279.28568
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 640.71432
$ws.Range("N107").Value = -5040

$ws.Range("H110").Value = 18000
$ws.Range("J110").Value = 18000
$ws.Range("L110").Value = 18000
$ws.Range("N110").Value = -26180

$ws.Range("H113").Value = 95553
$ws.Range("I113").Value = 71581.71000000001
$ws.Range("J113").Value = 143495.58
$ws.Range("K113").Value = 214745.13
$ws.Range("L113").Value = 430486.74
$ws.Range("M113").Value = -212575.13
$ws.Range("N113").Value = -434826.74

$ws.Range("H122").Value = 9526520
$ws.Range("I122").Value = 14287657
$ws.Range("J122").Value = 4245.7144
$ws.Range("K122").Value = 42862971
$ws.Range("L122").Value = 12737.1432
$ws.Range("M122").Value = -42860521
$ws.Range("N122").Value = -17637.1432

Write-Output "Applied 35 row updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
